$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 14.1720012268567
$ws.Cells.Item(2, 3).Value = 8.275842683704788
$ws.Cells.Item(2, 4).Value = 14.52841793321709
$ws.Cells.Item(2, 5).Value = 15.72943841080991
$ws.Cells.Item(2, 7).Value = 31.37917305893154
$ws.Cells.Item(2, 8).Value = 14.95655497564942
$ws.Cells.Item(2, 10).Value = 9.243027905734245
$ws.Cells.Item(2, 13).Value = 17.56313531573104
$ws.Cells.Item(2, 15).Value = 23.11630886565293
$ws.Cells.Item(3, 2).Value = 13.5834780133244
$ws.Cells.Item(3, 3).Value = 7.781640683800221
$ws.Cells.Item(3, 4).Value = 14.52137030850117
$ws.Cells.Item(3, 5).Value = 15.75997132645237
$ws.Cells.Item(3, 7).Value = 31.4900940657927
$ws.Cells.Item(3, 8).Value = 15.02130489557099
$ws.Cells.Item(3, 10).Value = 9.26897574624601
$ws.Cells.Item(3, 13).Value = 17.3697679801714
$ws.Cells.Item(3, 15).Value = 23.22161923524536
$ws.Cells.Item(4, 2).Value = 13.20948019419935
$ws.Cells.Item(4, 3).Value = 7.460540841914348
$ws.Cells.Item(4, 4).Value = 14.52017427488449
$ws.Cells.Item(4, 5).Value = 15.7818265098357
$ws.Cells.Item(4, 7).Value = 31.57116720568767
$ws.Cells.Item(4, 8).Value = 15.06416389685467
$ws.Cells.Item(4, 10).Value = 9.285872907385082
$ws.Cells.Item(4, 13).Value = 17.25214261610362
$ws.Cells.Item(4, 15).Value = 23.29271503194428
$ws.Cells.Item(5, 2).Value = 13.05409819687535
$ws.Cells.Item(5, 3).Value = 7.325290210596325
$ws.Cells.Item(5, 4).Value = 14.52047561049705
$ws.Cells.Item(5, 5).Value = 15.79151312478434
$ws.Cells.Item(5, 7).Value = 31.60744360370208
$ws.Cells.Item(5, 8).Value = 15.08240839995464
$ws.Cells.Item(5, 10).Value = 9.293001787516937
$ws.Cells.Item(5, 13).Value = 17.20453210823743
$ws.Cells.Item(5, 15).Value = 23.32329916653133
$ws.Cells.Item(6, 2).Value = 13.02812394192755
$ws.Cells.Item(6, 3).Value = 7.302567611665298
$ws.Cells.Item(6, 4).Value = 14.52057331453197
$ws.Cells.Item(6, 5).Value = 15.79316868372949
$ws.Cells.Item(6, 7).Value = 31.61366218685114
$ws.Cells.Item(6, 8).Value = 15.08548491373884
$ws.Cells.Item(6, 10).Value = 9.294200232226189
$ws.Cells.Item(6, 13).Value = 17.19664723212742
$ws.Cells.Item(6, 15).Value = 23.32847482190336
$ws.Cells.Item(7, 2).Value = 13.20739641434802
$ws.Cells.Item(7, 3).Value = 7.458734557819917
$ws.Cells.Item(7, 4).Value = 14.52017514383409
$ws.Cells.Item(7, 5).Value = 15.78195398838975
$ws.Cells.Item(7, 7).Value = 31.57164336011274
$ws.Cells.Item(7, 8).Value = 15.0644067953184
$ws.Cells.Item(7, 10).Value = 9.285968064825489
$ws.Cells.Item(7, 13).Value = 17.25149915759078
$ws.Cells.Item(7, 15).Value = 23.29312098138158
$ws.Cells.Item(8, 2).Value = 13.97181682715822
$ws.Cells.Item(8, 3).Value = 8.10912758803388
$ws.Cells.Item(8, 4).Value = 14.52533897519714
$ws.Cells.Item(8, 5).Value = 15.73932053962108
$ws.Cells.Item(8, 7).Value = 31.41471432548978
$ws.Cells.Item(8, 8).Value = 14.97823627613217
$ws.Cells.Item(8, 10).Value = 9.251774702387214
$ws.Cells.Item(8, 13).Value = 17.49626073671856
$ws.Cells.Item(8, 15).Value = 23.15128071974177
$ws.Cells.Item(9, 2).Value = 15.3628327029193
$ws.Cells.Item(9, 3).Value = 9.243178334117417
$ws.Cells.Item(9, 4).Value = 14.56023151356515
$ws.Cells.Item(9, 5).Value = 15.6804176622967
$ws.Cells.Item(9, 7).Value = 31.21079078126259
$ws.Cells.Item(9, 8).Value = 14.83391544468147
$ws.Cells.Item(9, 10).Value = 9.192357688896461
$ws.Cells.Item(9, 13).Value = 17.98297678175942
$ws.Cells.Item(9, 15).Value = 22.92445689648408
$ws.Cells.Item(10, 2).Value = 16.31025900064291
$ws.Cells.Item(10, 3).Value = 9.989060742839028
$ws.Cells.Item(10, 4).Value = 14.60083155965695
$ws.Cells.Item(10, 5).Value = 15.65225000368977
$ws.Cells.Item(10, 7).Value = 31.12541021345467
$ws.Cells.Item(10, 8).Value = 14.74297626402155
$ws.Cells.Item(10, 10).Value = 9.153329336424767
$ws.Cells.Item(10, 13).Value = 18.34197599798113
$ws.Cells.Item(10, 15).Value = 22.7894576726909
$ws.Cells.Item(11, 2).Value = 16.72364505020793
$ws.Cells.Item(11, 3).Value = 10.3093098556072
$ws.Cells.Item(11, 4).Value = 14.62251225696514
$ws.Cells.Item(11, 5).Value = 15.6427238007096
$ws.Cells.Item(11, 7).Value = 31.10077081576105
$ws.Cells.Item(11, 8).Value = 14.70489632858358
$ws.Cells.Item(11, 10).Value = 9.136572415615717
$ws.Cells.Item(11, 13).Value = 18.50503134485161
$ws.Cells.Item(11, 15).Value = 22.73498991452633
$ws.Cells.Item(12, 2).Value = 16.87755567761319
$ws.Cells.Item(12, 3).Value = 10.42783504833598
$ws.Cells.Item(12, 4).Value = 14.63117985493701
$ws.Cells.Item(12, 5).Value = 15.63958947951148
$ws.Cells.Item(12, 7).Value = 31.09349666439295
$ws.Cells.Item(12, 8).Value = 14.69095054834205
$ws.Cells.Item(12, 10).Value = 9.130369952665648
$ws.Cells.Item(12, 13).Value = 18.56669258477967
$ws.Cells.Item(12, 15).Value = 22.71536893233182
$ws.Cells.Item(13, 2).Value = 16.84452663563855
$ws.Cells.Item(13, 3).Value = 10.40243072292921
$ws.Cells.Item(13, 4).Value = 14.62929285385334
$ws.Cells.Item(13, 5).Value = 15.64024346976322
$ws.Cells.Item(13, 7).Value = 31.09497162687714
$ws.Cells.Item(13, 8).Value = 14.69393290702906
$ws.Cells.Item(13, 10).Value = 9.131699409306854
$ws.Cells.Item(13, 13).Value = 18.55341722961755
$ws.Cells.Item(13, 15).Value = 22.71954987676436
$ws.Cells.Item(14, 2).Value = 16.73636058386495
$ws.Cells.Item(14, 3).Value = 10.31911605402478
$ws.Cells.Item(14, 4).Value = 14.62321620069052
$ws.Cells.Item(14, 5).Value = 15.64245645652563
$ws.Cells.Item(14, 7).Value = 31.1001310826887
$ws.Cells.Item(14, 8).Value = 14.7037394881279
$ws.Cells.Item(14, 10).Value = 9.136059271129335
$ws.Cells.Item(14, 13).Value = 18.5101062105395
$ws.Cells.Item(14, 15).Value = 22.73335551220439
$ws.Cells.Item(15, 2).Value = 16.66976062482206
$ws.Cells.Item(15, 3).Value = 10.26772566172737
$ws.Cells.Item(15, 4).Value = 14.6195535327537
$ws.Cells.Item(15, 5).Value = 15.64387358554286
$ws.Cells.Item(15, 7).Value = 31.10355956900263
$ws.Cells.Item(15, 8).Value = 14.7098081056286
$ws.Cells.Item(15, 10).Value = 9.138748426345275
$ws.Cells.Item(15, 13).Value = 18.4835645883816
$ws.Cells.Item(15, 15).Value = 22.741942896128
$ws.Cells.Item(16, 2).Value = 16.282880273344
$ws.Cells.Item(16, 3).Value = 9.967747660055954
$ws.Cells.Item(16, 4).Value = 14.59947894951751
$ws.Cells.Item(16, 5).Value = 15.65293874612464
$ws.Cells.Item(16, 7).Value = 31.12730747297823
$ws.Cells.Item(16, 8).Value = 14.74553117030475
$ws.Cells.Item(16, 10).Value = 9.154444475009738
$ws.Cells.Item(16, 13).Value = 18.33131072436046
$ws.Cells.Item(16, 15).Value = 22.79315753703516
$ws.Cells.Item(17, 2).Value = 16.04095825423234
$ws.Cells.Item(17, 3).Value = 9.778834703330833
$ws.Cells.Item(17, 4).Value = 14.5879835458196
$ws.Cells.Item(17, 5).Value = 15.65934219929609
$ws.Cells.Item(17, 7).Value = 31.14552395264338
$ws.Cells.Item(17, 8).Value = 14.76828939639345
$ws.Cells.Item(17, 10).Value = 9.16432864004685
$ws.Cells.Item(17, 13).Value = 18.23780841622799
$ws.Cells.Item(17, 15).Value = 22.8263590820545
$ws.Cells.Item(18, 2).Value = 15.90015908375217
$ws.Cells.Item(18, 3).Value = 9.66838337610473
$ws.Cells.Item(18, 4).Value = 14.58167426324006
$ws.Cells.Item(18, 5).Value = 15.66333471504098
$ws.Cells.Item(18, 7).Value = 31.15733750254931
$ws.Cells.Item(18, 8).Value = 14.78168879789953
$ws.Cells.Item(18, 10).Value = 9.17010764104168
$ws.Cells.Item(18, 13).Value = 18.18400713838987
$ws.Cells.Item(18, 15).Value = 22.8461089777871
$ws.Cells.Item(19, 2).Value = 15.85220641718026
$ws.Cells.Item(19, 3).Value = 9.630678641300911
$ws.Cells.Item(19, 4).Value = 14.57959013574701
$ws.Cells.Item(19, 5).Value = 15.66473964244063
$ws.Cells.Item(19, 7).Value = 31.1615663442394
$ws.Cells.Item(19, 8).Value = 14.78627871629042
$ws.Cells.Item(19, 10).Value = 9.172080449222793
$ws.Cells.Item(19, 13).Value = 18.16578876934175
$ws.Cells.Item(19, 15).Value = 22.85290795638659
$ws.Cells.Item(20, 2).Value = 16.06688300531036
$ws.Cells.Item(20, 3).Value = 9.799130493643263
$ws.Cells.Item(20, 4).Value = 14.58917596537811
$ws.Cells.Item(20, 5).Value = 15.65862851416735
$ws.Cells.Item(20, 7).Value = 31.14344640559471
$ws.Cells.Item(20, 8).Value = 14.76583470547037
$ws.Cells.Item(20, 10).Value = 9.163266739294491
$ws.Cells.Item(20, 13).Value = 18.24776444061599
$ws.Cells.Item(20, 15).Value = 22.8227570668434
$ws.Cells.Item(21, 2).Value = 16.76820365479866
$ws.Cells.Item(21, 3).Value = 10.34366215285264
$ws.Cells.Item(21, 4).Value = 14.62498867835683
$ws.Cells.Item(21, 5).Value = 15.64179360874416
$ws.Cells.Item(21, 7).Value = 31.09855971654946
$ws.Cells.Item(21, 8).Value = 14.70084617534976
$ws.Cells.Item(21, 10).Value = 9.134774795575645
$ws.Cells.Item(21, 13).Value = 18.52283036869062
$ws.Cells.Item(21, 15).Value = 22.72927314425085
$ws.Cells.Item(22, 2).Value = 17.21119387080197
$ws.Cells.Item(22, 3).Value = 10.68353767910769
$ws.Cells.Item(22, 4).Value = 14.65105922207674
$ws.Cells.Item(22, 5).Value = 15.63354830846863
$ws.Cells.Item(22, 7).Value = 31.08121416611068
$ws.Cells.Item(22, 8).Value = 14.66113742469393
$ws.Cells.Item(22, 10).Value = 9.116987114552806
$ws.Cells.Item(22, 13).Value = 18.70208922337925
$ws.Cells.Item(22, 15).Value = 22.6740353730383
$ws.Cells.Item(23, 2).Value = 16.97619613898544
$ws.Cells.Item(23, 3).Value = 10.5036054849805
$ws.Cells.Item(23, 4).Value = 14.63690257130625
$ws.Cells.Item(23, 5).Value = 15.63769663189487
$ws.Cells.Item(23, 7).Value = 31.08937053851333
$ws.Cells.Item(23, 8).Value = 14.68207728611629
$ws.Cells.Item(23, 10).Value = 9.126404599929314
$ws.Cells.Item(23, 13).Value = 18.60647777544419
$ws.Cells.Item(23, 15).Value = 22.70297866195288
$ws.Cells.Item(24, 2).Value = 16.05516775951233
$ws.Cells.Item(24, 3).Value = 9.789960506187356
$ws.Cells.Item(24, 4).Value = 14.58863593902963
$ws.Cells.Item(24, 5).Value = 15.65895020223435
$ws.Cells.Item(24, 7).Value = 31.14438148962172
$ws.Cells.Item(24, 8).Value = 14.76694348887379
$ws.Cells.Item(24, 10).Value = 9.163746524113312
$ws.Cells.Item(24, 13).Value = 18.24326346024303
$ws.Cells.Item(24, 15).Value = 22.82438347661558
$ws.Cells.Item(25, 2).Value = 14.99907033689091
$ws.Cells.Item(25, 3).Value = 8.951663461500655
$ws.Cells.Item(25, 4).Value = 14.54815256744984
$ws.Cells.Item(25, 5).Value = 15.69370221869306
$ws.Cells.Item(25, 7).Value = 31.25471745697357
$ws.Cells.Item(25, 8).Value = 14.87031192619341
$ws.Cells.Item(25, 10).Value = 9.207617180919911
$ws.Cells.Item(25, 13).Value = 17.85086493789654
$ws.Cells.Item(25, 15).Value = 22.98028655481615
